$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies A1:D18, with row 1 being headers and rows 2-18 being data.
# The edit sorts the data rows (A2:D18) in ascending order by column A (time),
# keeping each row's B/C/D values together with its A value.
$dataRange = $ws.Range("A2:D18")
$values = $dataRange.Value2

$rowCount = $values.GetLength(0)
$colCount = $values.GetLength(1)

$rows = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $rowCount; $r++) {
    $row = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $row += $values[$r, $c]
    }
    $rows.Add($row) | Out-Null
}

$sorted = $rows | Sort-Object { [double]$_[0] }

$out = New-Object 'object[,]' $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $out[$r, $c] = $sorted[$r][$c]
    }
}

$dataRange.Value2 = $out
